$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each "block" of curl rows has its own header cell in column B (e.g. B6, B40,
# B74, B107, B140) computed from CONCATENATE("https://.../wu-map/", B<mapid>, "/").
# The per-row helper cells in column B within a block should just reference
# that block's own header cell (relative reference), but several of them were
# incorrectly hard-coded to the very first block's header ($B$6), which made
# every block after the first show the wrong map id. Fix them to reference
# their own block header with a plain relative formula.

$blocks = @(
    @{ Header = "B6";   Rows = @(7,8,9,12,16,17,18,19,20,21,22,23,27,28,29) },
    @{ Header = "B40";  Rows = @(41,42,43,46,50,51,52,53,54,55,56,57,61,62,63) },
    @{ Header = "B74";  Rows = @(75,76,77,80,84,85,86,87,88,89,90,91,95,96,97) },
    @{ Header = "B107"; Rows = @(108,109,110,113,117,118,119,120,121,122,123,124,128,129,130) },
    @{ Header = "B140"; Rows = @(141,142,143,146,150,151,152,153,154,155,156,157,161,162,163) }
)

foreach ($block in $blocks) {
    $headerRef = $block.Header
    foreach ($row in $block.Rows) {
        $ws.Range("B$row").Formula = "=$headerRef"
    }
}

# Update the saved selection state (last cell the editor happened to land on).
$ws.Range("A141").Select() | Out-Null
